# 7.10 Fixed Some Bugs
# Clear the stray debug values that were left in H4/I4 (a 0 placeholder and
# the "Dee-Thinking1" avatar tag) and move the active selection onto that
# range, matching where the bug actually was.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H4:I4").ClearContents()

$ws.Range("H4:I4").Select()
